$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows 2-5 (existing rows get new values) and new row 6 appended.
$data = @(
    @(1, 8, 2, 4, 4, -4, 2, 23, 5),
    @(2, 7, 2, 2, 3, -5, 1, 12, 5),
    @(3, 6, 3, 5, 8, -1, 5, 56, 5),
    @(4, 9, 4, 7, 8, -2, 4, 45, 5),
    @(5, 7, 0, 4, 3, -3, 3, 34, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
    # Column J (10) holds the "version" shared string "train_dim2_1"
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

# Update the selection shown in the saved sheet view.
$ws.Range("I1").Select() | Out-Null
